$d = $word.ActiveDocument

# Replace the TODO placeholder bullet with the new assumption about the
# U-Haul (or other rental company) web API.
$d.Content.Find.Execute(
    "TODO: Flush out assumptions",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The U-Haul company or the other rental companies have a web api that we can " + [char]8220 + "plug" + [char]8221 + " into for data, otherwise we have to create a scraper",
    2
)
